# Multi-browser support: extend the TestsRunner sheet with Browser/Execute
# columns plus per-browser credential columns, add two more test rows,
# widen the new Browser column, and update the recorded cell selections.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("TestsRunner")

# --- TestsRunner (sheet2): header row -------------------------------------
$ws2.Range("A1").Value = "TestCase"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "Browser"
$ws2.Range("D1").Value = "Execute"
$ws2.Range("E1").Value = "username"
$ws2.Range("F1").Value = "password"

# --- TestsRunner (sheet2): data rows ---------------------------------------
$ws2.Range("A2").Value = "loginLogoutTest"
$ws2.Range("B2").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C2").Value = "chrome"
$ws2.Range("D2").Value = "no"
$ws2.Range("E2").Value = "admin"
$ws2.Range("F2").Value = "ad123"

$ws2.Range("A3").Value = "loginLogoutTest"
$ws2.Range("B3").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C3").Value = "chrome"
$ws2.Range("D3").Value = "no"
$ws2.Range("E3").Value = "ad123"
$ws2.Range("F3").Value = "admin"

$ws2.Range("A4").Value = "loginLogoutTest"
$ws2.Range("B4").Value = "validate OrangeHRM login and logout functionality"
$ws2.Range("C4").Value = "firefox"
$ws2.Range("D4").Value = "yes"
$ws2.Range("E4").Value = "Admin"
$ws2.Range("F4").Value = "admin123"

$ws2.Range("A5").Value = "homePageTitleTest"
$ws2.Range("B5").Value = "validate title of home page"
$ws2.Range("C5").Value = "chrome"
$ws2.Range("D5").Value = "yes"
$ws2.Range("E5").Value = "Admin"
$ws2.Range("F5").Value = "admin123"

# New "Browser" column is a bit wider than the default.
$ws2.Columns.Item(3).ColumnWidth = 15.33

# --- Selections -------------------------------------------------------------
# Sheet1's remembered selection moved one cell to the right ...
$ws1.Range("C8").Select()
# ... and TestsRunner (the active tab) ends up selected at D11.
$ws2.Range("D11").Select()
